# Sealing_Stainless_Steel_Tamper-Resistant_Button_Head_Torx_Screws.xlsx
#
# The sheet currently has a header row (row 1, bold/bordered style) with
# descriptive labels ("Length", "Threading", ...) in A1:O1, followed by the
# data rows.
#
# The target layout instead keeps row 1 as a numeric index row (0..14, same
# bold/bordered style) and pushes the old descriptive header text down into
# a brand-new row 2 (plain style, and with the last two header labels -
# "thread_size"/"material_surface" - dropped rather than carried over).
# Every row that used to start at row 2 shifts down by one to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current row-1 header text (Length, Threading, ...) before we
# overwrite row 1 - we need it again for the new row 2.
$headers = @()
for ($col = 1; $col -le 15; $col++) {
    $headers += $ws.Cells.Item(1, $col).Value2
}

# Insert a fresh row at position 2. This shifts the old rows 2-25 (material
# note, size rows, ...) down to rows 3-26, carrying their content/format
# with them, and grows the sheet from O25 to O26.
$ws.Rows.Item(2).Insert()

# The inserted row inherits row 1's (bold/bordered) formatting by default;
# strip that back to the plain/default style used by the rest of the data
# rows.
$ws.Rows.Item(2).ClearFormats()

# Row 1 becomes a plain numeric index: 0, 1, 2, ... 14 (keeps style s=1).
for ($col = 1; $col -le 15; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# Row 2 gets the old header text in columns A-K and M (Pkg.). Columns L, N
# and O (Pkg.Qty. part number column, thread_size, material_surface) are
# left blank in the new row 2.
$skipCols = @(12, 14, 15)
for ($col = 1; $col -le 13; $col++) {
    if ($skipCols -contains $col) { continue }
    $ws.Cells.Item(2, $col).Value = $headers[$col - 1]
}
